$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7491.5884
$ws.Range("J17").Value = 7491.5884
$ws.Range("L17").Value = 22474.7652
$ws.Range("N17").Value = -22810.7652
$ws.Range("H69").Value = 7679.7144
$ws.Range("J69").Value = 7843.0557
$ws.Range("L69").Value = 23529.1671
$ws.Range("N69").Value = -25277.1671
$ws.Range("H72").Value = 7679.7144
$ws.Range("J72").Value = 7843.0557
$ws.Range("L72").Value = 70587.5013
$ws.Range("N72").Value = -79323.5013
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 5022.2144
$ws.Range("J100").Value = 5702.5557
$ws.Range("L100").Value = 5702.5557
$ws.Range("N100").Value = -6784.5557
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H111").Value = 786.4286
$ws.Range("I111").Value = 709.1667
$ws.Range("K111").Value = 2127.5001
$ws.Range("M111").Value = 939.4998999999998
$ws.Range("H112").Value = 2822.3635
$ws.Range("J112").Value = 2822.3635
$ws.Range("L112").Value = 8467.0905
$ws.Range("N112").Value = -10683.0905
$ws.Range("H138").Value = 4076.6885
$ws.Range("I138").Value = 2870
$ws.Range("J138").Value = 4258.83
$ws.Range("K138").Value = 8610
$ws.Range("L138").Value = 12776.49
$ws.Range("M138").Value = -3470
$ws.Range("N138").Value = -23056.49
$ws.Range("H141").Value = 2143.3684
$ws.Range("I141").Value = 1302.6428
$ws.Range("K141").Value = 3907.9284
$ws.Range("M141").Value = 1272.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 7300.3413
$ws.Range("I32").Value = 7300.3413
$ws.Range("K32").Value = 7300.3413
$ws.Range("M32").Value = -7013.3413
$ws.Range("H61").Value = 7400.9
$ws.Range("I61").Value = 6862.9375
$ws.Range("J61").Value = 9552.75
$ws.Range("K61").Value = 6862.9375
$ws.Range("L61").Value = 9552.75
$ws.Range("M61").Value = -6650.9375
$ws.Range("N61").Value = -9976.75
$ws.Range("H132").Value = 4985.615
$ws.Range("I132").Value = 3816.5833
$ws.Range("K132").Value = 11449.7499
$ws.Range("M132").Value = -8919.749899999999
$ws.Range("H136").Value = 7400.9
$ws.Range("I136").Value = 6862.9375
$ws.Range("J136").Value = 9552.75
$ws.Range("K136").Value = 20588.8125
$ws.Range("L136").Value = 28658.25
$ws.Range("M136").Value = -18038.8125
$ws.Range("N136").Value = -33758.25
$ws.Range("H139").Value = 68513.28999999999
$ws.Range("J139").Value = 68513.28999999999
$ws.Range("L139").Value = 68513.28999999999
$ws.Range("N139").Value = -78793.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 7999.5
$ws.Range("I5").Value = 8000
$ws.Range("J5").Value = 7999
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 7999
$ws.Range("M5").Value = -7887
$ws.Range("N5").Value = -8225
$ws.Range("H95").Value = 44666.332
$ws.Range("J95").Value = 44666.332
$ws.Range("L95").Value = 44666.332
$ws.Range("N95").Value = -50158.332
$ws.Range("H134").Value = 3762.2354
$ws.Range("I134").Value = 3026.5386
$ws.Range("J134").Value = 6153.25
$ws.Range("K134").Value = 9079.6158
$ws.Range("L134").Value = 18459.75
$ws.Range("M134").Value = -6544.6158
$ws.Range("N134").Value = -23529.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38716.633
$ws.Range("I31").Value = 2702.9285
$ws.Range("J31").Value = 70228.625
$ws.Range("K31").Value = 2702.9285
$ws.Range("L31").Value = 70228.625
$ws.Range("M31").Value = -2407.9285
$ws.Range("N31").Value = -70818.625
$ws.Range("H34").Value = 38716.633
$ws.Range("I34").Value = 2702.9285
$ws.Range("J34").Value = 70228.625
$ws.Range("K34").Value = 2702.9285
$ws.Range("L34").Value = 70228.625
$ws.Range("M34").Value = -2500.9285
$ws.Range("N34").Value = -70632.625
$ws.Range("H58").Value = 4931.2
$ws.Range("I58").Value = 3344.9092
$ws.Range("J58").Value = 7615.6924
$ws.Range("K58").Value = 3344.9092
$ws.Range("L58").Value = 7615.6924
$ws.Range("M58").Value = -3141.9092
$ws.Range("N58").Value = -8021.6924
$ws.Range("H133").Value = 55234.734
$ws.Range("J133").Value = 55293.418
$ws.Range("L133").Value = 55293.418
$ws.Range("N133").Value = -60353.418
$ws.Range("H134").Value = 3008.7273
$ws.Range("I134").Value = 2123.4119
$ws.Range("K134").Value = 6370.2357
$ws.Range("M134").Value = -3835.2357
$ws.Range("H136").Value = 4931.2
$ws.Range("I136").Value = 3344.9092
$ws.Range("J136").Value = 7615.6924
$ws.Range("K136").Value = 10034.7276
$ws.Range("L136").Value = 22847.0772
$ws.Range("M136").Value = -7484.7276
$ws.Range("N136").Value = -27947.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1013.5
$ws.Range("I110").Value = 1013.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3040.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1049.5
$ws.Range("N110").ClearContents()
$ws.Range("H124").Value = 6730.8335
$ws.Range("J124").Value = 6877
$ws.Range("L124").Value = 20631
$ws.Range("N124").Value = -30451
$ws.Range("H137").Value = 62275.766
$ws.Range("I137").Value = 1874.875
$ws.Range("J137").Value = 115965.445
$ws.Range("K137").Value = 5624.625
$ws.Range("L137").Value = 347896.335
$ws.Range("M137").Value = -524.625
$ws.Range("N137").Value = -358096.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H12").Value = 22000000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H122").Value = 5031.696
$ws.Range("I122").Value = 5041.846
$ws.Range("K122").Value = 15125.538
$ws.Range("M122").Value = -12675.538
$ws.Range("H126").Value = 4068.08
$ws.Range("I126").Value = 2957.077
$ws.Range("J126").Value = 5271.6665
$ws.Range("K126").Value = 8871.231
$ws.Range("L126").Value = 15814.9995
$ws.Range("M126").Value = -6401.231
$ws.Range("N126").Value = -20754.9995
$ws.Range("H132").Value = 7415.3335
$ws.Range("I132").Value = 4566.647
$ws.Range("K132").Value = 13699.941
$ws.Range("M132").Value = -11169.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5899.026
$ws.Range("I7").Value = 4487.8335
$ws.Range("K7").Value = 4487.8335
$ws.Range("M7").Value = -4375.8335
$ws.Range("H22").Value = 3585.8386
$ws.Range("I22").Value = 1180.8667
$ws.Range("K22").Value = 1180.8667
$ws.Range("M22").Value = -885.8667
$ws.Range("H27").Value = 3585.8386
$ws.Range("I27").Value = 1180.8667
$ws.Range("K27").Value = 1180.8667
$ws.Range("M27").Value = -1073.8667
$ws.Range("H46").Value = 4309.1816
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4309.1816
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4309.1816
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4685.1816
$ws.Range("H61").Value = 12730.7
$ws.Range("I61").Value = 8550.333000000001
$ws.Range("J61").Value = 19001.25
$ws.Range("K61").Value = 8550.333000000001
$ws.Range("L61").Value = 19001.25
$ws.Range("M61").Value = -8348.333000000001
$ws.Range("N61").Value = -19405.25
$ws.Range("H108").Value = 73311
$ws.Range("J108").Value = 84967
$ws.Range("L108").Value = 84967
$ws.Range("N108").Value = -92647
$ws.Range("H112").Value = 52631
$ws.Range("J112").Value = 52631
$ws.Range("L112").Value = 52631
$ws.Range("N112").Value = -55585
$ws.Range("H113").Value = 12730.7
$ws.Range("I113").Value = 8550.333000000001
$ws.Range("J113").Value = 19001.25
$ws.Range("K113").Value = 8550.333000000001
$ws.Range("L113").Value = 19001.25
$ws.Range("M113").Value = -6380.333000000001
$ws.Range("N113").Value = -23341.25
$ws.Range("H114").Value = 52631
$ws.Range("J114").Value = 52631
$ws.Range("L114").Value = 52631
$ws.Range("N114").Value = -61309
$ws.Range("H126").Value = 5899.026
$ws.Range("I126").Value = 4487.8335
$ws.Range("K126").Value = 13463.5005
$ws.Range("M126").Value = -10993.5005
$ws.Range("H136").Value = 4319.875
$ws.Range("I136").Value = 3165.1667
$ws.Range("J136").Value = 7784
$ws.Range("K136").Value = 9495.500100000001
$ws.Range("L136").Value = 23352
$ws.Range("M136").Value = -6945.500100000001
$ws.Range("N136").Value = -28452

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1068.75
$ws.Range("I107").Value = 874.2593000000001
$ws.Range("K107").Value = 2622.7779
$ws.Range("M107").Value = -702.7779
